$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the old placeholder content in A1
$ws.Range("A1").ClearContents()

# Populate the match table (order chosen to match the shared-string table order)
$ws.Range("E4").Value = "Prabhu"
$ws.Range("G4").Value = "Mani"
$ws.Range("H4").Value = "Sendil"
$ws.Range("I4").Value = "Thamim"
$ws.Range("J4").Value = "Dinesh"
$ws.Range("K4").Value = "Ranjit"
$ws.Range("L4").Value = "vadivel"
$ws.Range("M4").Value = "Ganapathi"
$ws.Range("N4").Value = "prabha"
$ws.Range("O4").Value = "Abi"
$ws.Range("F4").Value = "Venu"
$ws.Range("C4").Value = "Total pts"
$ws.Range("D4").Value = "games"
$ws.Range("C5").Value = "19(1,0)"
$ws.Range("C6").Value = "25(1,0)"
$ws.Range("I7").Value = "25(6)"
$ws.Range("J5").Value = "19(8)"
$ws.Range("E10").Value = "20(8)"
$ws.Range("M6").Value = "25(7)"
$ws.Range("F13").Value = "2(7)"
$ws.Range("C13").Value = "2(0,1)"
$ws.Range("C9").Value = "4(0,1)"
$ws.Range("G9").Value = "4(6)"
$ws.Range("C3").Value = "(win,los)"
$ws.Range("B5").Value = "Prabhu"
$ws.Range("B6").Value = "Venu"
$ws.Range("B7").Value = "Mani"
$ws.Range("C7").Value = "25(1,0)"
$ws.Range("B8").Value = "Sendil"
$ws.Range("B9").Value = "Thamim"
$ws.Range("B10").Value = "Dinesh"
$ws.Range("B11").Value = "Ranjit"
$ws.Range("B12").Value = "vadivel"
$ws.Range("B13").Value = "Ganapathi"
$ws.Range("B14").Value = "prabha"
$ws.Range("B15").Value = "Abi"

# Apply a thin box border around every cell in the table range
$tableRange = $ws.Range("B3:O15")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# Select the table range, matching the saved selection state
$tableRange.Select() | Out-Null
